# Auto-generated edit script: refresh market-price derived columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 6798.8887
$ws.Range("I28").Value = 4189.75
$ws.Range("J28").Value = 8886.200000000001
$ws.Range("K28").Value = 4189.75
$ws.Range("L28").Value = 8886.200000000001
$ws.Range("M28").Value = -3704.75
$ws.Range("N28").Value = -9856.200000000001

$ws.Range("H94").Value = 3432.7778
$ws.Range("I94").Value = 2986.875
$ws.Range("J94").Value = 7000
$ws.Range("K94").Value = 2986.875
$ws.Range("L94").Value = 7000
$ws.Range("M94").Value = -2535.875

$ws.Range("H107").Value = 1235.3077
$ws.Range("I107").Value = 579.0909
$ws.Range("J107").Value = 4844.5
$ws.Range("K107").Value = 579.0909
$ws.Range("L107").Value = 4844.5
$ws.Range("M107").Value = 1340.9091
$ws.Range("N107").Value = -8684.5

$ws.Range("H111").Value = 7391
$ws.Range("I111").Value = 1682.3334
$ws.Range("J111").Value = 12284.143
$ws.Range("K111").Value = 5047.0002
$ws.Range("L111").Value = 36852.429
$ws.Range("M111").Value = -1980.0002
$ws.Range("N111").Value = -42986.429

$ws.Range("H115").Value = 1374.75
$ws.Range("I115").Value = 750
$ws.Range("J115").Value = 1999.5
$ws.Range("K115").Value = 2250
$ws.Range("L115").Value = 5998.5
$ws.Range("M115").Value = -683
$ws.Range("N115").Value = -9132.5

$ws.Range("H132").Value = 4712.617
$ws.Range("I132").Value = 2477.6099
$ws.Range("J132").Value = 19985.166
$ws.Range("K132").Value = 7432.8297
$ws.Range("L132").Value = 59955.49800000001
$ws.Range("M132").Value = -4902.8297
$ws.Range("N132").Value = -65015.49800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""

$ws.Range("H110").Value = 4074.257
$ws.Range("I110").Value = 3930.2334
$ws.Range("J110").Value = 4938.4
$ws.Range("K110").Value = 3930.2334
$ws.Range("L110").Value = 4938.4
$ws.Range("M110").Value = -1885.2334

$ws.Range("H132").Value = 3504.0715
$ws.Range("I132").Value = 3547.6155
$ws.Range("J132").Value = 2938
$ws.Range("K132").Value = 10642.8465
$ws.Range("L132").Value = 8814
$ws.Range("M132").Value = -8112.8465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 64443.5
$ws.Range("I26").Value = 28888
$ws.Range("J26").Value = 99999
$ws.Range("K26").Value = 28888
$ws.Range("L26").Value = 99999
$ws.Range("M26").Value = -28596

$ws.Range("H105").Value = 489400.34
$ws.Range("I105").Value = 881491.6
$ws.Range("J105").Value = 3954
$ws.Range("K105").Value = 881491.6
$ws.Range("L105").Value = 3954
$ws.Range("M105").Value = -879744.6
$ws.Range("N105").Value = -7448

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 365.33334
$ws.Range("I7").Value = 286
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 286
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -173
$ws.Range("N7").Value = -1226

$ws.Range("H31").Value = 15386750
$ws.Range("I31").Value = 23811684
$ws.Range("J31").Value = 2084.913
$ws.Range("K31").Value = 23811684
$ws.Range("L31").Value = 2084.913
$ws.Range("M31").Value = -23811389
$ws.Range("N31").Value = -2674.913

$ws.Range("H34").Value = 15386750
$ws.Range("I34").Value = 23811684
$ws.Range("J34").Value = 2084.913
$ws.Range("K34").Value = 23811684
$ws.Range("L34").Value = 2084.913
$ws.Range("M34").Value = -23811482
$ws.Range("N34").Value = -2488.913

$ws.Range("H99").Value = 24095.908
$ws.Range("I99").Value = 9149.5
$ws.Range("J99").Value = 32636.715
$ws.Range("K99").Value = 9149.5
$ws.Range("L99").Value = 32636.715
$ws.Range("M99").Value = -7651.5
$ws.Range("N99").Value = -35632.715

$ws.Range("H126").Value = 24095.908
$ws.Range("I126").Value = 9149.5
$ws.Range("J126").Value = 32636.715
$ws.Range("K126").Value = 27448.5
$ws.Range("L126").Value = 97910.145
$ws.Range("M126").Value = -24978.5
$ws.Range("N126").Value = -102850.145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 15089.388
$ws.Range("I56").Value = 15089.388
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 15089.388
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -14559.388

$ws.Range("H97").Value = 1093.0322
$ws.Range("I97").Value = 1570.75
$ws.Range("J97").Value = 1022.2593
$ws.Range("K97").Value = 4712.25
$ws.Range("L97").Value = 3066.7779
$ws.Range("M97").Value = -4216.25
$ws.Range("N97").Value = -4058.7779

$ws.Range("H123").Value = 12266.4
$ws.Range("I123").Value = 3999.5
$ws.Range("J123").Value = 17777.666
$ws.Range("K123").Value = 11998.5
$ws.Range("L123").Value = 53332.99800000001
$ws.Range("M123").Value = -9548.5

$ws.Range("H129").Value = 17192224
$ws.Range("I129").Value = 22729886
$ws.Range("J129").Value = 5009366.5
$ws.Range("K129").Value = 68189658
$ws.Range("L129").Value = 15028099.5
$ws.Range("M129").Value = -68184658

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 39999
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 39999
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 39999
$ws.Range("N44").Value = -41191

$ws.Range("H132").Value = 2397.2903
$ws.Range("I132").Value = 2250.4644
$ws.Range("J132").Value = 3767.6667
$ws.Range("K132").Value = 6751.3932
$ws.Range("L132").Value = 11303.0001
$ws.Range("M132").Value = -4221.3932
$ws.Range("N132").Value = -16363.0001

$ws.Range("H134").Value = 65303.57
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 65303.57
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 195910.71
$ws.Range("N134").Value = -200980.71

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3483.1738
$ws.Range("I122").Value = 2839.5
$ws.Range("J122").Value = 7774.3335
$ws.Range("K122").Value = 8518.5
$ws.Range("L122").Value = 23323.0005
$ws.Range("M122").Value = -6068.5
$ws.Range("N122").Value = -28223.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 20000
$ws.Range("I64").Value = 20000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 20000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -19752
$ws.Range("N64").Value = ""

$ws.Range("H67").Value = 20000
$ws.Range("I67").Value = 20000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 20000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -19142
$ws.Range("N67").Value = ""

$ws.Range("H100").Value = 970.1429000000001
$ws.Range("I100").Value = 965.1667
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 1930.3334
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1389.3334
$ws.Range("N100").Value = -3082

$ws.Range("H113").Value = 592.7778
$ws.Range("I113").Value = 433.14285
$ws.Range("J113").Value = 1151.5
$ws.Range("K113").Value = 1299.42855
$ws.Range("L113").Value = 3454.5
$ws.Range("M113").Value = 870.5714499999999
$ws.Range("N113").Value = -7794.5

$ws.Range("H132").Value = 1912.415
$ws.Range("I132").Value = 1379.9048
$ws.Range("J132").Value = 3945.6365
$ws.Range("K132").Value = 4139.7144
$ws.Range("L132").Value = 11836.9095
$ws.Range("M132").Value = -1609.7144
$ws.Range("N132").Value = -16896.9095
